$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'41.400.44"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.05%  "

# Row 3
$ws.Range("D3").Value = "'2.194.65"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'252.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.28%  "

# Row 6
$ws.Range("D6").Value = "'0.626"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.20%  "

# Row 7
$ws.Range("D7").Value = "'68.83"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.20%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.586"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.12%  "

# Row 10
$ws.Range("D10").Value = "'37.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.29%  "

# Row 11
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'58.23"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.26%  "

# Row 12
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "'0.0947"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.24%  "

# Row 13
$ws.Range("D13").Value = "'7.17"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.04%  "

# Row 15
$ws.Range("D15").Value = "'2.520.68"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.14%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'14.73"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.880"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.63%  "

# Row 18
$ws.Range("D18").Value = "'2.202.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.07%  "

# Row 19
$ws.Range("D19").Value = "'41.357.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.13%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0952"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.48%  "

# Row 21
$ws.Range("D21").Value = "'6.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.17%  "

# Row 22
$ws.Range("D22").Value = "'71.85"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.01%  "

# Row 23
$ws.Range("D23").Value = "'232.91"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.90%  "

# Row 24
$ws.Range("E24").Value = "  +3.07%  "

# Row 25
$ws.Range("D25").Value = "'12.12"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +22.43%  "

# Row 26
$ws.Range("D26").Value = "'3.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.33%  "

# Row 27
$ws.Range("E27").Value = "  +0.05%  "

# Row 28
$ws.Range("D28").Value = "'2.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.34%  "

# Row 29
$ws.Range("E29").Value = "  -2.48%  "

# Row 30
$ws.Range("D30").Value = "'170.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.37%  "

# Row 31
$ws.Range("D31").Value = "'20.67"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.24%  "

# Row 32
$ws.Range("E32").Value = "  +1.04%  "

# Row 33
$ws.Range("D33").Value = "'5.57"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.42%  "

# Row 34
$ws.Range("D34").Value = "'0.123"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.86%  "

# Row 35
$ws.Range("D35").Value = "'0.0730"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.90%  "

# Row 36
$ws.Range("D36").Value = "'26.39"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +15.19%  "

# Row 37
$ws.Range("D37").Value = "'4.62"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.36%  "

# Row 38
$ws.Range("D38").Value = "'4.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.16%  "

# Row 39
$ws.Range("D39").Value = "'0.0299"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.48%  "

# Row 40
$ws.Range("E40").Value = "  -2.06%  "

# Row 41
$ws.Range("D41").Value = "'5.78"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.85%  "

# Row 42
$ws.Range("D42").Value = "'12.06"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +18.56%  "

# Row 43
$ws.Range("D43").Value = "'64.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.59%  "

# Row 44
$ws.Range("D44").Value = "'5.02"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.64%  "

# Row 45
$ws.Range("D45").Value = "'0.204"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.54%  "

# Row 46
$ws.Range("D46").Value = "'8.65"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.28%  "

# Row 47
$ws.Range("E47").Value = "  +0.50%  "

# Row 48
$ws.Range("D48").Value = "'1.01"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.64%  "

# Row 49
$ws.Range("D49").Value = "'1.15"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.18%  "

# Row 50
$ws.Range("D50").Value = "'1.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.70%  "

# Row 51
$ws.Range("D51").Value = "'2.32"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.17%  "
